$wb = $excel.ActiveWorkbook
$games = $wb.Worksheets.Item("Games")
$next = $wb.Worksheets.Item("Next")

# Season results updated through 1/17: the game that was scheduled next
# (DET on 2024-01-15, serial 45306) has now been played. Append it as a
# completed game (row 40) on the "Games" sheet, and remove it from the
# "Next" sheet (row 2), shifting the remaining scheduled games up by one row.

$newRow = 40
$games.Cells.Item($newRow, 1).Value = 39
$dateCell = $games.Cells.Item($newRow, 2)
$dateCell.Value = 45306
$dateCell.NumberFormat = "YYYY-MM-DD"
$games.Cells.Item($newRow, 3).Value = -1
$games.Cells.Item($newRow, 4).Value = 117
$games.Cells.Item($newRow, 5).Value = 99.3
$games.Cells.Item($newRow, 6).Value = 0.5629999999999999
$games.Cells.Item($newRow, 7).Value = 7.4
$games.Cells.Item($newRow, 8).Value = 16.7
$games.Cells.Item($newRow, 9).Value = 0.105
$games.Cells.Item($newRow, 10).Value = 117.8
$games.Cells.Item($newRow, 11).Value = "DET"
$games.Cells.Item($newRow, 12).Value = 129
$games.Cells.Item($newRow, 13).Value = 0.622
$games.Cells.Item($newRow, 14).Value = 13
$games.Cells.Item($newRow, 15).Value = 32.4
$games.Cells.Item($newRow, 16).Value = 0.189
$games.Cells.Item($newRow, 17).Value = 129.9
$games.Cells.Item($newRow, 18).Value = 1
$games.Cells.Item($newRow, 19).Value = 0

# Remove the now-played game from the "Next" sheet, shifting remaining rows up.
$next.Rows.Item(2).Delete()
